$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new Slides topics (7.0 Power, 6.1 / 6.2 Logical Effort) ---
$ws.Range("C13").Value = "7.0 - Power"
$ws.Range("C13").HorizontalAlignment = -4131  # xlLeft

$ws.Range("C11").Value = "6.1 - Logical Effort Part 1"
$ws.Range("C11").HorizontalAlignment = -4131  # xlLeft

$ws.Range("C12").Value = "6.2 - Logical Effort Part 2"
$ws.Range("C12").HorizontalAlignment = -4131  # xlLeft

# --- schedule / due-date cleanup ---
# "Lab 1" due date moved from Links(G13) to Due(E13)
$ws.Range("G13").ClearContents()

# "Lab  1" (double space) due under Links for week 14
$ws.Range("G15").Value = "Lab  1"

# "Quiz 1" due under Links for week 13 (was combined into "HW 1/Quiz 1" on G12)
$ws.Range("G14").Value = "Quiz 1"

$ws.Range("E13").Value = "Lab 1"

# Week 11 Due column now just "HW 1" (quiz moved out)
$ws.Range("G12").Value = "HW 1"

# --- new Zoom recording links ---
$ws.Hyperlinks.Add($ws.Range("H11"), "https://iu.zoom.us/rec/share/IB0cyZD_eISgemVTiFcyU12VHnYR-nnvtp9ufdfrSWPA2uJRqo_G2Z5NIcCkVlR9.C4kofYzOqIlKnWWL")
$ws.Range("H11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("H12"), "https://iu.zoom.us/rec/share/3GW27LneiHNfTu4meorvl4ZkINGOM493rWDRJw5Tc2lGb7ikOhb-Y3GSUW2vEk-n.Kp6h406uHpjtf0FM")
$ws.Range("H12").Style = "Hyperlink"

# --- restore cursor/selection position ---
$ws.Range("H23").Select() | Out-Null
